$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (Unmet Demand / Wasted Surplus / Household Surplus) ---
$ws.Range("D2").Value = 0.08408200188452446
$ws.Range("E2").Value = 0.7100534011061919
$ws.Range("F2").Value = 49701999.65585799

$ws.Range("D3").Value = 0.08417940453896246
$ws.Range("E3").Value = 0.7100534011061919
$ws.Range("F3").Value = 49695750.08687391

$ws.Range("D4").Value = 0.2182006796012734
$ws.Range("E4").Value = 0.7100534011061919
$ws.Range("F4").Value = 42140680.16732705

$ws.Range("D5").Value = 0.4785207646748159
$ws.Range("E5").Value = 0.7100534011061919
$ws.Range("F5").Value = 27302166.61425338

# --- Add new rows 6, 7, 8 ---
$ws.Range("A6").Value = "40"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0.46
$ws.Range("D6").Value = 0.6088905735061112
$ws.Range("E6").Value = 0.7100534011061919
$ws.Range("F6").Value = 14454088.20754593

$ws.Range("A7").Value = "50"
$ws.Range("A8").Value = "60"

# --- Styling ---
# Currency style for Household Surplus column
$ws.Range("F2:F6").Style = "Currency"

# Percent style for Unmet Demand / Wasted Surplus columns
$ws.Range("D2:E6").Style = "Percent"

# RE target column (A) keeps its header-like bold/border/center formatting,
# with percentage number formatting layered on top
$ws.Range("A2:A8").Style = "Percent"
$ws.Range("A2:A8").Font.Bold = $true
$ws.Range("A2:A8").Borders.LineStyle = 1
$ws.Range("A2:A8").HorizontalAlignment = -4108
$ws.Range("A2:A8").VerticalAlignment = -4160
$ws.Range("A2:A8").NumberFormat = "0%"

# --- Column widths (auto-fit like Excel does after entering data) ---
$ws.Columns("D:F").AutoFit()

# --- Selection ---
$ws.Range("K10").Select()

Write-Output "done"
